$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap interleaved duplicate-fixture rows (home/away got shuffled on re-scrape) ---
# rows 19 and 20
$ws.Cells.Item(19, 6).Value = "Barito Putera"
$ws.Cells.Item(19, 7).Value = 3
$ws.Cells.Item(19, 8).Value = "PSS Sleman"
$ws.Cells.Item(19, 9).Value = 1
$ws.Cells.Item(19, 10).Value = 1.56
$ws.Cells.Item(19, 11).Value = "12/07/2023 22:12"
$ws.Cells.Item(19, 12).Value = 1.86
$ws.Cells.Item(19, 13).Value = "14/07/2023 09:58"
$ws.Cells.Item(19, 14).Value = 4.07
$ws.Cells.Item(19, 15).Value = "12/07/2023 22:12"
$ws.Cells.Item(19, 16).Value = 3.27
$ws.Cells.Item(19, 17).Value = "14/07/2023 09:57"
$ws.Cells.Item(19, 18).Value = 4.9
$ws.Cells.Item(19, 19).Value = "12/07/2023 22:12"
$ws.Cells.Item(19, 20).Value = 3.05
$ws.Cells.Item(19, 21).Value = "14/07/2023 09:58"
$ws.Cells.Item(19, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/ps-barito-putera-pss-sleman/nZin4VzL/"
$ws.Cells.Item(20, 6).Value = "Persikabo 1973"
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = "PSM Makassar"
$ws.Cells.Item(20, 9).Value = 1
$ws.Cells.Item(20, 10).Value = 3.73
$ws.Cells.Item(20, 11).Value = "12/07/2023 22:12"
$ws.Cells.Item(20, 12).Value = 3.53
$ws.Cells.Item(20, 13).Value = "14/07/2023 09:53"
$ws.Cells.Item(20, 14).Value = 3.49
$ws.Cells.Item(20, 15).Value = "12/07/2023 22:12"
$ws.Cells.Item(20, 16).Value = 3.33
$ws.Cells.Item(20, 17).Value = "14/07/2023 09:53"
$ws.Cells.Item(20, 18).Value = 1.86
$ws.Cells.Item(20, 19).Value = "12/07/2023 22:12"
$ws.Cells.Item(20, 20).Value = 2.11
$ws.Cells.Item(20, 21).Value = "14/07/2023 09:53"
$ws.Cells.Item(20, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persikabo-1973-psm-makassar/QPjj3kjR/"

# rows 42 and 43
$ws.Cells.Item(42, 6).Value = "Persis Solo"
$ws.Cells.Item(42, 7).Value = 1
$ws.Cells.Item(42, 8).Value = "Arema FC"
$ws.Cells.Item(42, 9).Value = 1
$ws.Cells.Item(42, 10).Value = 1.65
$ws.Cells.Item(42, 11).Value = "28/07/2023 22:12"
$ws.Cells.Item(42, 12).Value = 1.61
$ws.Cells.Item(42, 13).Value = "30/07/2023 09:55"
$ws.Cells.Item(42, 14).Value = 3.8
$ws.Cells.Item(42, 15).Value = "28/07/2023 22:12"
$ws.Cells.Item(42, 16).Value = 4.15
$ws.Cells.Item(42, 17).Value = "30/07/2023 09:55"
$ws.Cells.Item(42, 18).Value = 4.48
$ws.Cells.Item(42, 19).Value = "28/07/2023 22:12"
$ws.Cells.Item(42, 20).Value = 5.06
$ws.Cells.Item(42, 21).Value = "30/07/2023 09:51"
$ws.Cells.Item(42, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persis-solo-arema-fc/tvOgLNBC/"
$ws.Cells.Item(43, 6).Value = "RANS Nusantara"
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = "PSS Sleman"
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 2.32
$ws.Cells.Item(43, 11).Value = "28/07/2023 22:12"
$ws.Cells.Item(43, 12).Value = 2.79
$ws.Cells.Item(43, 13).Value = "30/07/2023 09:55"
$ws.Cells.Item(43, 14).Value = 3.26
$ws.Cells.Item(43, 15).Value = "28/07/2023 22:12"
$ws.Cells.Item(43, 16).Value = 3.5
$ws.Cells.Item(43, 17).Value = "30/07/2023 09:58"
$ws.Cells.Item(43, 18).Value = 2.75
$ws.Cells.Item(43, 19).Value = "28/07/2023 22:12"
$ws.Cells.Item(43, 20).Value = 2.42
$ws.Cells.Item(43, 21).Value = "30/07/2023 09:55"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/rans-nusantara-pss-sleman/pCUpNqs0/"

# rows 74 and 75
$ws.Cells.Item(74, 6).Value = "FC Bhayangkara"
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = "Borneo"
$ws.Cells.Item(74, 9).Value = 2
$ws.Cells.Item(74, 10).Value = 2.81
$ws.Cells.Item(74, 11).Value = "17/08/2023 02:12"
$ws.Cells.Item(74, 12).Value = 3.19
$ws.Cells.Item(74, 13).Value = "18/08/2023 13:59"
$ws.Cells.Item(74, 14).Value = 3.28
$ws.Cells.Item(74, 15).Value = "17/08/2023 02:12"
$ws.Cells.Item(74, 16).Value = 3.29
$ws.Cells.Item(74, 17).Value = "18/08/2023 13:59"
$ws.Cells.Item(74, 18).Value = 2.32
$ws.Cells.Item(74, 19).Value = "17/08/2023 02:12"
$ws.Cells.Item(74, 20).Value = 2.28
$ws.Cells.Item(74, 21).Value = "18/08/2023 13:59"
$ws.Cells.Item(74, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/fc-bhayangkara-borneo/dC1KPcKe/"
$ws.Cells.Item(75, 6).Value = "Persita"
$ws.Cells.Item(75, 7).Value = 2
$ws.Cells.Item(75, 8).Value = "PSS Sleman"
$ws.Cells.Item(75, 9).Value = 3
$ws.Cells.Item(75, 10).Value = 1.75
$ws.Cells.Item(75, 11).Value = "17/08/2023 02:12"
$ws.Cells.Item(75, 12).Value = 1.91
$ws.Cells.Item(75, 13).Value = "18/08/2023 13:51"
$ws.Cells.Item(75, 14).Value = 3.54
$ws.Cells.Item(75, 15).Value = "17/08/2023 02:12"
$ws.Cells.Item(75, 16).Value = 3.35
$ws.Cells.Item(75, 17).Value = "18/08/2023 13:51"
$ws.Cells.Item(75, 18).Value = 4.19
$ws.Cells.Item(75, 19).Value = "17/08/2023 02:12"
$ws.Cells.Item(75, 20).Value = 4.25
$ws.Cells.Item(75, 21).Value = "18/08/2023 13:51"
$ws.Cells.Item(75, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persita-pss-sleman/EmCFQw5k/"

# rows 82 and 83
$ws.Cells.Item(82, 6).Value = "Madura United"
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = "FC Bhayangkara"
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 1.76
$ws.Cells.Item(82, 11).Value = "23/08/2023 22:12"
$ws.Cells.Item(82, 12).Value = 1.59
$ws.Cells.Item(82, 13).Value = "25/08/2023 09:55"
$ws.Cells.Item(82, 14).Value = 3.57
$ws.Cells.Item(82, 15).Value = "23/08/2023 22:12"
$ws.Cells.Item(82, 16).Value = 4
$ws.Cells.Item(82, 17).Value = "25/08/2023 09:55"
$ws.Cells.Item(82, 18).Value = 3.95
$ws.Cells.Item(82, 19).Value = "23/08/2023 22:12"
$ws.Cells.Item(82, 20).Value = 5.59
$ws.Cells.Item(82, 21).Value = "25/08/2023 09:55"
$ws.Cells.Item(82, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/madura-united-fc-bhayangkara/AFRgvcZl/"
$ws.Cells.Item(83, 6).Value = "Persik Kediri"
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = "PSIS Semarang"
$ws.Cells.Item(83, 9).Value = 1
$ws.Cells.Item(83, 10).Value = 2.11
$ws.Cells.Item(83, 11).Value = "23/08/2023 22:12"
$ws.Cells.Item(83, 12).Value = 1.86
$ws.Cells.Item(83, 13).Value = "25/08/2023 09:58"
$ws.Cells.Item(83, 14).Value = 3.28
$ws.Cells.Item(83, 15).Value = "23/08/2023 22:12"
$ws.Cells.Item(83, 16).Value = 3.69
$ws.Cells.Item(83, 17).Value = "25/08/2023 09:58"
$ws.Cells.Item(83, 18).Value = 3.11
$ws.Cells.Item(83, 19).Value = "23/08/2023 22:12"
$ws.Cells.Item(83, 20).Value = 4.02
$ws.Cells.Item(83, 21).Value = "25/08/2023 09:58"
$ws.Cells.Item(83, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persik-kediri-psis-semarang/SjxkuwKr/"

# rows 84 and 85
$ws.Cells.Item(84, 6).Value = "Dewa United"
$ws.Cells.Item(84, 7).Value = 2
$ws.Cells.Item(84, 8).Value = "Persija Jakarta"
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 3.11
$ws.Cells.Item(84, 11).Value = "24/08/2023 02:12"
$ws.Cells.Item(84, 12).Value = 2.6
$ws.Cells.Item(84, 13).Value = "25/08/2023 13:59"
$ws.Cells.Item(84, 14).Value = 3.17
$ws.Cells.Item(84, 15).Value = "24/08/2023 02:12"
$ws.Cells.Item(84, 16).Value = 3.11
$ws.Cells.Item(84, 17).Value = "25/08/2023 13:59"
$ws.Cells.Item(84, 18).Value = 2.15
$ws.Cells.Item(84, 19).Value = "24/08/2023 02:12"
$ws.Cells.Item(84, 20).Value = 2.84
$ws.Cells.Item(84, 21).Value = "25/08/2023 13:59"
$ws.Cells.Item(84, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/dewa-united-persija-jakarta/IuV1xy41/"
$ws.Cells.Item(85, 6).Value = "Borneo"
$ws.Cells.Item(85, 7).Value = 2
$ws.Cells.Item(85, 8).Value = "Persita"
$ws.Cells.Item(85, 9).Value = 1
$ws.Cells.Item(85, 10).Value = 1.85
$ws.Cells.Item(85, 11).Value = "24/08/2023 02:12"
$ws.Cells.Item(85, 12).Value = 1.64
$ws.Cells.Item(85, 13).Value = "25/08/2023 13:51"
$ws.Cells.Item(85, 14).Value = 3.73
$ws.Cells.Item(85, 15).Value = "24/08/2023 02:12"
$ws.Cells.Item(85, 16).Value = 3.9
$ws.Cells.Item(85, 17).Value = "25/08/2023 13:51"
$ws.Cells.Item(85, 18).Value = 3.51
$ws.Cells.Item(85, 19).Value = "24/08/2023 02:12"
$ws.Cells.Item(85, 20).Value = 5.27
$ws.Cells.Item(85, 21).Value = "25/08/2023 13:51"
$ws.Cells.Item(85, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/borneo-persita/l6QcwHle/"

# rows 127 and 128
$ws.Cells.Item(127, 6).Value = "RANS Nusantara"
$ws.Cells.Item(127, 7).Value = 2
$ws.Cells.Item(127, 8).Value = "PSIS Semarang"
$ws.Cells.Item(127, 9).Value = 1
$ws.Cells.Item(127, 10).Value = 2.54
$ws.Cells.Item(127, 11).Value = "04/10/2023 21:12"
$ws.Cells.Item(127, 12).Value = 3.07
$ws.Cells.Item(127, 13).Value = "06/10/2023 09:52"
$ws.Cells.Item(127, 14).Value = 3.19
$ws.Cells.Item(127, 15).Value = "04/10/2023 21:12"
$ws.Cells.Item(127, 16).Value = 3.19
$ws.Cells.Item(127, 17).Value = "06/10/2023 09:52"
$ws.Cells.Item(127, 18).Value = 2.54
$ws.Cells.Item(127, 19).Value = "04/10/2023 21:12"
$ws.Cells.Item(127, 20).Value = 2.39
$ws.Cells.Item(127, 21).Value = "06/10/2023 09:52"
$ws.Cells.Item(127, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/rans-nusantara-psis-semarang/j15nYO7i/"
$ws.Cells.Item(128, 6).Value = "Persikabo 1973"
$ws.Cells.Item(128, 7).Value = 2
$ws.Cells.Item(128, 8).Value = "Persis Solo"
$ws.Cells.Item(128, 9).Value = 2
$ws.Cells.Item(128, 10).Value = 2.54
$ws.Cells.Item(128, 11).Value = "04/10/2023 21:12"
$ws.Cells.Item(128, 12).Value = 3.74
$ws.Cells.Item(128, 13).Value = "06/10/2023 09:56"
$ws.Cells.Item(128, 14).Value = 3.24
$ws.Cells.Item(128, 15).Value = "04/10/2023 21:12"
$ws.Cells.Item(128, 16).Value = 3.73
$ws.Cells.Item(128, 17).Value = "06/10/2023 09:58"
$ws.Cells.Item(128, 18).Value = 2.52
$ws.Cells.Item(128, 19).Value = "04/10/2023 21:12"
$ws.Cells.Item(128, 20).Value = 1.92
$ws.Cells.Item(128, 21).Value = "06/10/2023 09:56"
$ws.Cells.Item(128, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persikabo-1973-persis-solo/OE3fW2x4/"

# rows 139 and 140
$ws.Cells.Item(139, 6).Value = "FC Bhayangkara"
$ws.Cells.Item(139, 7).Value = 1
$ws.Cells.Item(139, 8).Value = "Barito Putera"
$ws.Cells.Item(139, 9).Value = 1
$ws.Cells.Item(139, 10).Value = 2.53
$ws.Cells.Item(139, 11).Value = "19/10/2023 21:12"
$ws.Cells.Item(139, 12).Value = 2.9
$ws.Cells.Item(139, 13).Value = "21/10/2023 09:59"
$ws.Cells.Item(139, 14).Value = 3.21
$ws.Cells.Item(139, 15).Value = "19/10/2023 21:12"
$ws.Cells.Item(139, 16).Value = 3.22
$ws.Cells.Item(139, 17).Value = "21/10/2023 09:57"
$ws.Cells.Item(139, 18).Value = 2.53
$ws.Cells.Item(139, 19).Value = "19/10/2023 21:12"
$ws.Cells.Item(139, 20).Value = 2.49
$ws.Cells.Item(139, 21).Value = "21/10/2023 09:59"
$ws.Cells.Item(139, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/fc-bhayangkara-ps-barito-putera/80bvWgrL/"
$ws.Cells.Item(140, 6).Value = "PSS Sleman"
$ws.Cells.Item(140, 7).Value = 2
$ws.Cells.Item(140, 8).Value = "Persik Kediri"
$ws.Cells.Item(140, 9).Value = 2
$ws.Cells.Item(140, 10).Value = 2.52
$ws.Cells.Item(140, 11).Value = "19/10/2023 21:12"
$ws.Cells.Item(140, 12).Value = 2.28
$ws.Cells.Item(140, 13).Value = "21/10/2023 09:59"
$ws.Cells.Item(140, 14).Value = 3.09
$ws.Cells.Item(140, 15).Value = "19/10/2023 21:12"
$ws.Cells.Item(140, 16).Value = 3.29
$ws.Cells.Item(140, 17).Value = "21/10/2023 09:59"
$ws.Cells.Item(140, 18).Value = 2.63
$ws.Cells.Item(140, 19).Value = "19/10/2023 21:12"
$ws.Cells.Item(140, 20).Value = 3.17
$ws.Cells.Item(140, 21).Value = "21/10/2023 09:55"
$ws.Cells.Item(140, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/pss-sleman-persik-kediri/nicrVDcR/"

# --- Append two new match rows (167, 168), copying formatting from the last existing row ---
$ws.Range("A166:V166").Copy()
$ws.Range("A167:V167").PasteSpecial(-4122)
$ws.Cells.Item(167, 1).Value = 166
$ws.Cells.Item(167, 2).Value = "indonesia"
$ws.Cells.Item(167, 3).Value = "liga-1"
$ws.Cells.Item(167, 4).Value = "2023-2024"
$ws.Cells.Item(167, 5).Value = 45239.54166666666
$ws.Cells.Item(167, 6).Value = "Persija Jakarta"
$ws.Cells.Item(167, 7).Value = 4
$ws.Cells.Item(167, 8).Value = "Persikabo 1973"
$ws.Cells.Item(167, 9).Value = 0
$ws.Cells.Item(167, 10).Value = 1.38
$ws.Cells.Item(167, 11).Value = "08/11/2023 01:12"
$ws.Cells.Item(167, 12).Value = 1.49
$ws.Cells.Item(167, 13).Value = "09/11/2023 12:57"
$ws.Cells.Item(167, 14).Value = 4.55
$ws.Cells.Item(167, 15).Value = "08/11/2023 01:12"
$ws.Cells.Item(167, 16).Value = 4.39
$ws.Cells.Item(167, 17).Value = "09/11/2023 12:59"
$ws.Cells.Item(167, 18).Value = 6.25
$ws.Cells.Item(167, 19).Value = "08/11/2023 01:12"
$ws.Cells.Item(167, 20).Value = 6.39
$ws.Cells.Item(167, 21).Value = "09/11/2023 12:59"
$ws.Cells.Item(167, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persija-jakarta-persikabo-1973/YkI6s8dt/"

$ws.Range("A166:V166").Copy()
$ws.Range("A168:V168").PasteSpecial(-4122)
$ws.Cells.Item(168, 1).Value = 167
$ws.Cells.Item(168, 2).Value = "indonesia"
$ws.Cells.Item(168, 3).Value = "liga-1"
$ws.Cells.Item(168, 4).Value = "2023-2024"
$ws.Cells.Item(168, 5).Value = 45239.54166666666
$ws.Cells.Item(168, 6).Value = "RANS Nusantara"
$ws.Cells.Item(168, 7).Value = 1
$ws.Cells.Item(168, 8).Value = "FC Bhayangkara"
$ws.Cells.Item(168, 9).Value = 1
$ws.Cells.Item(168, 10).Value = 1.96
$ws.Cells.Item(168, 11).Value = "08/11/2023 01:12"
$ws.Cells.Item(168, 12).Value = 1.94
$ws.Cells.Item(168, 13).Value = "09/11/2023 12:58"
$ws.Cells.Item(168, 14).Value = 3.45
$ws.Cells.Item(168, 15).Value = "08/11/2023 01:12"
$ws.Cells.Item(168, 16).Value = 3.59
$ws.Cells.Item(168, 17).Value = "09/11/2023 12:54"
$ws.Cells.Item(168, 18).Value = 3.32
$ws.Cells.Item(168, 19).Value = "08/11/2023 01:12"
$ws.Cells.Item(168, 20).Value = 3.77
$ws.Cells.Item(168, 21).Value = "09/11/2023 12:58"
$ws.Cells.Item(168, 22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/rans-nusantara-fc-bhayangkara/tGLEunRh/"

$ws.Application.CutCopyMode = $false